$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(70960, 71021, 71031, 71037, 71038, 71057, 71084)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

[void]$ws.Range("A2:A8").Select()
